# Insert a new weekly data row for "Macroferia Regional de Talca - Choclo"
# right above the existing row 152, shifting the subsequent rows (old 152-214)
# down to (153-215) and filling the newly opened row 152 with the new week's
# observation (Fecha 44636 = 2022-03-16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 152..214 down by one, duplicating row 151's formatting into the
# freshly inserted row 152 (keeps the date-style on column D, etc.).
$ws.Rows("152").Insert()

# Populate the new row with this week's data.
$ws.Cells.Item(152, 1).Value  = 5
$ws.Cells.Item(152, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(152, 3).Value  = "Maule"
$ws.Cells.Item(152, 4).Value  = 44636
$ws.Cells.Item(152, 5).Value  = 7
$ws.Cells.Item(152, 6).Value  = 100112024
$ws.Cells.Item(152, 7).Value  = "Choclo"
$ws.Cells.Item(152, 8).Value  = "Choclero"
$ws.Cells.Item(152, 9).Value  = "Primera"
$ws.Cells.Item(152, 10).Value = 30000
$ws.Cells.Item(152, 11).Value = 120
$ws.Cells.Item(152, 12).Value = 120
$ws.Cells.Item(152, 13).Value = 120
$ws.Cells.Item(152, 14).Value = "$/unidad"
$ws.Cells.Item(152, 15).Value = "Región del Maule"
$ws.Cells.Item(152, 16).Value = 120
$ws.Cells.Item(152, 17).Value = 1
$ws.Cells.Item(152, 18).Value = "Hortaliza"
